$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 3): date, duration (time), and label
$ws.Range("A3").Value = 44316
$ws.Range("A3").NumberFormat = "d-mmm"

$ws.Range("B3").Value = 0.059027777777777783
$ws.Range("B3").NumberFormat = "h:mm"

$ws.Range("C3").Value = "mise en page statique"

# Update the active selection to mirror the author's final cursor position
$ws.Range("C4").Select()
